$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.173780487804878
$ws.Range("C2").Value = 0.5823170731707317
$ws.Range("J2").Value = 0.01219512195121951
$ws.Range("P2").Value = 0.1219512195121951
$ws.Range("S2").Value = 0.1097560975609756
$ws.Range("B3").Value = 0.0101010101010101
$ws.Range("C3").Value = 0.02525252525252525
$ws.Range("J3").Value = 0.04040404040404041
$ws.Range("P3").Value = 0.7373737373737373
$ws.Range("S3").Value = 0.1868686868686869
$ws.Range("J4").Value = 0.1219512195121951
$ws.Range("P4").Value = 0.8048780487804879
$ws.Range("S4").Value = 0.07317073170731707
$ws.Range("J5").Value = 0.25
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.06319702602230483
$ws.Range("D6").Value = 0.003717472118959108
$ws.Range("F6").Value = 0.09665427509293681
$ws.Range("J6").Value = 0.2527881040892193
$ws.Range("O6").Value = 0.03345724907063197
$ws.Range("Q6").Value = 0.09665427509293681
$ws.Range("R6").Value = 0.05576208178438662
$ws.Range("S6").Value = 0.3977695167286245
$ws.Range("B7").Value = 0.06787330316742081
$ws.Range("D7").Value = 0.01357466063348416
$ws.Range("F7").Value = 0.06787330316742081
$ws.Range("J7").Value = 0.1131221719457014
$ws.Range("O7").Value = 0.01809954751131222
$ws.Range("Q7").Value = 0.2307692307692308
$ws.Range("R7").Value = 0.09049773755656108
$ws.Range("S7").Value = 0.3981900452488688
$ws.Range("B8").Value = 0.1092278719397363
$ws.Range("D8").Value = 0.01694915254237288
$ws.Range("E8").Value = 0.003766478342749529
$ws.Range("F8").Value = 0.07344632768361582
$ws.Range("J8").Value = 0.1035781544256121
$ws.Range("O8").Value = 0.01318267419962335
$ws.Range("Q8").Value = 0.1826741996233522
$ws.Range("R8").Value = 0.07344632768361582
$ws.Range("S8").Value = 0.423728813559322
$ws.Range("B9").Value = 0.07075471698113207
$ws.Range("D9").Value = 0.02358490566037736
$ws.Range("E9").Value = 0.004716981132075472
$ws.Range("F9").Value = 0.0660377358490566
$ws.Range("J9").Value = 0.1509433962264151
$ws.Range("O9").Value = 0.02830188679245283
$ws.Range("Q9").Value = 0.1933962264150944
$ws.Range("R9").Value = 0.08962264150943396
$ws.Range("S9").Value = 0.3726415094339622
$ws.Range("B10").Value = 0.1079584775086505
$ws.Range("D10").Value = 0.01730103806228374
$ws.Range("E10").Value = 0.001384083044982699
$ws.Range("F10").Value = 0.06643598615916955
$ws.Range("J10").Value = 0.1397923875432526
$ws.Range("O10").Value = 0.01730103806228374
$ws.Range("Q10").Value = 0.2179930795847751
$ws.Range("R10").Value = 0.0754325259515571
$ws.Range("S10").Value = 0.356401384083045
$ws.Range("G11").Value = 0.152046783625731
$ws.Range("J11").Value = 0.07894736842105263
$ws.Range("K11").Value = 0.2046783625730994
$ws.Range("L11").Value = 0.5467836257309941
$ws.Range("S11").Value = 0.01754385964912281
$ws.Range("G12").Value = 0.7295918367346939
$ws.Range("J12").Value = 0.1887755102040816
$ws.Range("K12").Value = 0.02040816326530612
$ws.Range("L12").Value = 0.04591836734693878
$ws.Range("S12").Value = 0.01530612244897959
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2444444444444444
$ws.Range("S13").Value = 0.08888888888888889
$ws.Range("F15").Value = 0.02325581395348837
$ws.Range("H15").Value = 0.1279069767441861
$ws.Range("I15").Value = 0.06589147286821706
$ws.Range("J15").Value = 0.3372093023255814
$ws.Range("K15").Value = 0.07751937984496124
$ws.Range("M15").Value = 0.007751937984496124
$ws.Range("O15").Value = 0.1085271317829457
$ws.Range("S15").Value = 0.251937984496124
$ws.Range("F16").Value = 0.01886792452830189
$ws.Range("H16").Value = 0.2122641509433962
$ws.Range("I16").Value = 0.04716981132075472
$ws.Range("J16").Value = 0.4009433962264151
$ws.Range("K16").Value = 0.1084905660377359
$ws.Range("M16").Value = 0.01415094339622642
$ws.Range("N16").Value = 0.004716981132075472
$ws.Range("O16").Value = 0.0660377358490566
$ws.Range("S16").Value = 0.1273584905660377
$ws.Range("F17").Value = 0.01509433962264151
$ws.Range("H17").Value = 0.1547169811320755
$ws.Range("I17").Value = 0.0830188679245283
$ws.Range("J17").Value = 0.4339622641509434
$ws.Range("K17").Value = 0.1113207547169811
$ws.Range("M17").Value = 0.009433962264150943
$ws.Range("N17").Value = 0.001886792452830189
$ws.Range("O17").Value = 0.06792452830188679
$ws.Range("S17").Value = 0.1226415094339623
$ws.Range("F18").Value = 0.01507537688442211
$ws.Range("H18").Value = 0.185929648241206
$ws.Range("I18").Value = 0.07537688442211055
$ws.Range("J18").Value = 0.4723618090452261
$ws.Range("K18").Value = 0.07035175879396985
$ws.Range("M18").Value = 0.01005025125628141
$ws.Range("O18").Value = 0.04522613065326633
$ws.Range("S18").Value = 0.1256281407035176
$ws.Range("F19").Value = 0.02367688022284123
$ws.Range("H19").Value = 0.2325905292479109
$ws.Range("I19").Value = 0.08844011142061281
$ws.Range("J19").Value = 0.346100278551532
$ws.Range("K19").Value = 0.1044568245125348
$ws.Range("M19").Value = 0.02367688022284123
$ws.Range("O19").Value = 0.06406685236768803
$ws.Range("S19").Value = 0.116991643454039
